$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "MEC-3B-Usin. CNC"
$ws.Range("E2").Value = "-"

$ws.Range("B3").Value = "-"
$ws.Range("C3").Value = "MEC-3B-Usin. CNC"
$ws.Range("E3").Value = "-"

$ws.Range("B4").Value = "-"
$ws.Range("C4").Value = "-"
$ws.Range("D4").Value = "MEC-2B-Ajustagem"

$ws.Range("C6").Value = "-"
$ws.Range("D6").Value = "MEC-2B-Ajustagem"

$ws.Range("B7").Value = "MEC-3B-Usin. CNC"
$ws.Range("C7").Value = "-"
$ws.Range("D7").Value = "MEC-2B-Ajustagem"

$ws.Range("B8").Value = "MEC-3B-Usin. CNC"
$ws.Range("C8").Value = "-"
$ws.Range("D8").Value = "MEC-2B-Ajustagem"
